$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 93, shifting existing rows 93:184 down to 94:185
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with the new data record
$ws.Range("A93").Value = 8
$ws.Range("B93").Value = "Terminal La Palmera de La Serena"
$ws.Range("C93").Value = "Coquimbo"
$ws.Range("D93").Value = 44589
$ws.Range("E93").Value = 4
$ws.Range("F93").Value = 100112021
$ws.Range("G93").Value = "Ají"
$ws.Range("H93").Value = "Americana (o)"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 520
$ws.Range("K93").Value = 10000
$ws.Range("L93").Value = 11000
$ws.Range("M93").Value = 10500
$ws.Range("N93").Value = "$/caja 15 kilos"
$ws.Range("O93").Value = "Provincia de Limarí"
$ws.Range("P93").Value = 700
$ws.Range("Q93").Value = 15
$ws.Range("R93").Value = "Hortaliza"
